$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 9).Value = "sv"
$ws.Cells.Item(8, 10).Value = "Statement-opinion"
$ws.Cells.Item(15, 9).Value = "sv"
$ws.Cells.Item(15, 10).Value = "Statement-opinion"
$ws.Cells.Item(20, 9).Value = "sd"
$ws.Cells.Item(20, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(23, 9).Value = "sd"
$ws.Cells.Item(23, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(28, 9).Value = "aa"
$ws.Cells.Item(28, 10).Value = "Agree/Accept"
$ws.Cells.Item(36, 9).Value = "ba"
$ws.Cells.Item(36, 10).Value = "Appreciation"
$ws.Cells.Item(41, 9).Value = "sv"
$ws.Cells.Item(41, 10).Value = "Statement-opinion"
$ws.Cells.Item(55, 9).Value = "sd"
$ws.Cells.Item(55, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(57, 9).Value = "sv"
$ws.Cells.Item(57, 10).Value = "Statement-opinion"
$ws.Cells.Item(62, 9).Value = "sd"
$ws.Cells.Item(62, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(71, 9).Value = "sd"
$ws.Cells.Item(71, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(73, 9).Value = "sd"
$ws.Cells.Item(73, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(81, 9).Value = "sd"
$ws.Cells.Item(81, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(84, 9).Value = "sv"
$ws.Cells.Item(84, 10).Value = "Statement-opinion"
$ws.Cells.Item(92, 9).Value = "b"
$ws.Cells.Item(92, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(93, 9).Value = "b"
$ws.Cells.Item(93, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(105, 9).Value = "aa"
$ws.Cells.Item(105, 10).Value = "Agree/Accept"
$ws.Cells.Item(109, 9).Value = "sv"
$ws.Cells.Item(109, 10).Value = "Statement-opinion"
$ws.Cells.Item(118, 9).Value = "aa"
$ws.Cells.Item(118, 10).Value = "Agree/Accept"
$ws.Cells.Item(126, 9).Value = "aa"
$ws.Cells.Item(126, 10).Value = "Agree/Accept"
$ws.Cells.Item(129, 9).Value = "sd"
$ws.Cells.Item(129, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(130, 9).Value = "sd"
$ws.Cells.Item(130, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(133, 9).Value = "sd"
$ws.Cells.Item(133, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(134, 9).Value = "%"
$ws.Cells.Item(134, 10).Value = "Uninterpretable"
$ws.Cells.Item(151, 9).Value = "sd"
$ws.Cells.Item(151, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(169, 9).Value = "sv"
$ws.Cells.Item(169, 10).Value = "Statement-opinion"
$ws.Cells.Item(173, 9).Value = "sv"
$ws.Cells.Item(173, 10).Value = "Statement-opinion"
$ws.Cells.Item(179, 9).Value = "sd"
$ws.Cells.Item(179, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(212, 9).Value = "sd"
$ws.Cells.Item(212, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(213, 9).Value = "sd"
$ws.Cells.Item(213, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(215, 9).Value = "sv"
$ws.Cells.Item(215, 10).Value = "Statement-opinion"
$ws.Cells.Item(216, 9).Value = "sd"
$ws.Cells.Item(216, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(225, 9).Value = "sv"
$ws.Cells.Item(225, 10).Value = "Statement-opinion"
$ws.Cells.Item(232, 9).Value = "aa"
$ws.Cells.Item(232, 10).Value = "Agree/Accept"
$ws.Cells.Item(236, 9).Value = "ba"
$ws.Cells.Item(236, 10).Value = "Appreciation"
$ws.Cells.Item(241, 9).Value = "aa"
$ws.Cells.Item(241, 10).Value = "Agree/Accept"
$ws.Cells.Item(242, 9).Value = "aa"
$ws.Cells.Item(242, 10).Value = "Agree/Accept"
$ws.Cells.Item(243, 9).Value = "aa"
$ws.Cells.Item(243, 10).Value = "Agree/Accept"
$ws.Cells.Item(249, 9).Value = "aa"
$ws.Cells.Item(249, 10).Value = "Agree/Accept"
$ws.Cells.Item(259, 9).Value = "aa"
$ws.Cells.Item(259, 10).Value = "Agree/Accept"
$ws.Cells.Item(264, 9).Value = "sd"
$ws.Cells.Item(264, 10).Value = "Statement-non-opinion"

$wb.Save()
